$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = 4000
$ws.Range("K16").Value = 4000
$ws.Range("M16").Value = -3770
$ws.Range("H31").Value = 4500
$ws.Range("I31").Value = 4500
$ws.Range("K31").Value = 13500
$ws.Range("M31").Value = -13270
$ws.Range("H33").Value = 496
$ws.Range("I33").Value = 281.4
$ws.Range("J33").Value = 1032.5
$ws.Range("K33").Value = 281.4
$ws.Range("L33").Value = 1032.5
$ws.Range("M33").Value = -52.39999999999998
$ws.Range("N33").Value = -1490.5
$ws.Range("H38").Value = 6928.875
$ws.Range("I38").Value = 6928.875
$ws.Range("K38").Value = 20786.625
$ws.Range("M38").Value = -20414.625
$ws.Range("H53").Value = 710.0909
$ws.Range("I53").Value = 959.4286
$ws.Range("K53").Value = 959.4286
$ws.Range("M53").Value = -322.4286
$ws.Range("H59").Value = 10000
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 30000
$ws.Range("M59").Value = -29443
$ws.Range("H80").Value = 1669.2778
$ws.Range("I80").Value = 1294.5555
$ws.Range("J80").Value = 2044
$ws.Range("K80").Value = 3883.6665
$ws.Range("L80").Value = 6132
$ws.Range("M80").Value = -2885.6665
$ws.Range("N80").Value = -8128
$ws.Range("H82").Value = 16666.666
$ws.Range("I82").Value = 16666.666
$ws.Range("K82").Value = 49999.99800000001
$ws.Range("M82").Value = -49593.99800000001
$ws.Range("H83").Value = 1669.2778
$ws.Range("I83").Value = 1294.5555
$ws.Range("J83").Value = 2044
$ws.Range("K83").Value = 11650.9995
$ws.Range("L83").Value = 18396
$ws.Range("M83").Value = -6658.9995
$ws.Range("N83").Value = -28380
$ws.Range("H85").Value = 16666.666
$ws.Range("I85").Value = 16666.666
$ws.Range("K85").Value = 49999.99800000001
$ws.Range("M85").Value = -48595.99800000001
$ws.Range("H129").Value = 1798.6154
$ws.Range("I129").Value = 580.1818
$ws.Range("K129").Value = 1740.5454
$ws.Range("M129").Value = 3259.4546
$ws.Range("H132").Value = 7751.4287
$ws.Range("I132").Value = 8393.333
$ws.Range("K132").Value = 25179.999
$ws.Range("M132").Value = -22649.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 66
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 85
$ws.Range("I5").Value = 85
$ws.Range("K5").Value = 85
$ws.Range("M5").Value = 27
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H21").Value = 5487.5
$ws.Range("I21").Value = 8475
$ws.Range("J21").Value = 2500
$ws.Range("K21").Value = 8475
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = -8101
$ws.Range("N21").Value = -3248
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H35").Value = 2020.25
$ws.Range("I35").Value = 2020.25
$ws.Range("K35").Value = 2020.25
$ws.Range("M35").Value = -1614.25
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 85
$ws.Range("I4").Value = 85
$ws.Range("K4").Value = 85
$ws.Range("M4").Value = 30
$ws.Range("H16").Value = 1308.1
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 1219.0526
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 1219.0526
$ws.Range("M16").Value = -2830
$ws.Range("N16").Value = -1559.0526
$ws.Range("H22").Value = 270.5
$ws.Range("I22").Value = 270.5
$ws.Range("K22").Value = 270.5
$ws.Range("M22").Value = -97.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1000
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1226
$ws.Range("H10").Value = 5759.625
$ws.Range("I10").Value = 235.4
$ws.Range("K10").Value = 235.4
$ws.Range("M10").Value = -96.4
$ws.Range("H11").Value = 226.66667
$ws.Range("I11").Value = 90
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 90
$ws.Range("L11").Value = 500
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = -780
$ws.Range("H31").Value = 4041
$ws.Range("I31").Value = 3822.25
$ws.Range("J31").Value = 4332.6665
$ws.Range("K31").Value = 3822.25
$ws.Range("L31").Value = 4332.6665
$ws.Range("M31").Value = -3527.25
$ws.Range("N31").Value = -4922.6665
$ws.Range("H34").Value = 4041
$ws.Range("I34").Value = 3822.25
$ws.Range("J34").Value = 4332.6665
$ws.Range("K34").Value = 3822.25
$ws.Range("L34").Value = 4332.6665
$ws.Range("M34").Value = -3620.25
$ws.Range("N34").Value = -4736.6665
$ws.Range("H99").Value = 1299.6666
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H107").Value = 1801.8334
$ws.Range("H126").Value = 1299.6666
$ws.Range("J126").Value = 1500
$ws.Range("L126").Value = 4500
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12.6
$ws.Range("I2").Value = 11
$ws.Range("K2").Value = 66
$ws.Range("M2").Value = 47
$ws.Range("H5").Value = 864.5
$ws.Range("J5").Value = 914.6667
$ws.Range("L5").Value = 2744.0001
$ws.Range("N5").Value = -2968.0001
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = -2773
$ws.Range("H131").Value = 1655.909
$ws.Range("I131").Value = 943.3333
$ws.Range("J131").Value = 1923.125
$ws.Range("K131").Value = 2829.9999
$ws.Range("L131").Value = 5769.375
$ws.Range("M131").Value = 2210.0001
$ws.Range("N131").Value = -15849.375
$ws.Range("H135").Value = 864.5
$ws.Range("J135").Value = 914.6667
$ws.Range("L135").Value = 8232.0003
$ws.Range("N135").Value = -13302.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H69").Value = 90000
$ws.Range("I69").Value = 90000
$ws.Range("K69").Value = 90000
$ws.Range("M69").Value = -89251
$ws.Range("H72").Value = 90000
$ws.Range("I72").Value = 90000
$ws.Range("K72").Value = 270000
$ws.Range("M72").Value = -266256

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4997.5
$ws.Range("I93").Value = 4997.5
$ws.Range("K93").Value = 4997.5
$ws.Range("M93").Value = -3749.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 2324.75
$ws.Range("I132").Value = 1433
$ws.Range("K132").Value = 4299
$ws.Range("M132").Value = -1769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20003334
$ws.Range("J5").Value = 20003334
$ws.Range("L5").Value = 20003334
$ws.Range("N5").Value = -20003558
$ws.Range("H23").Value = 12020
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15458
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H132").Value = 2865.3333
$ws.Range("I132").Value = 2865.3333
$ws.Range("K132").Value = 8595.999899999999
$ws.Range("M132").Value = -6065.999899999999
